# Update trial-2 (row 3) values for the "red_square_00 / bottom_left" schedule:
#   y_corrSteps (E3): 5 -> 6
#   y_nrSteps   (G3): 2 -> 3
#   alienID     (H3): 14 -> 13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 13

# Move/leave the active selection on E3, matching the saved cursor position.
$ws.Range("E3").Select()
